# Remove the trailing "Ver no Jupiter Salvar em pdf Salvar em docx" block:
# an empty paragraph, the "Ver no Jupiter..." paragraph, another empty
# paragraph, and a page-break paragraph, which used to sit right after the
# "LOM3081: Introdução à Mecânica dos Sólidos (Requisito fraco)" requirement
# line.

$d = $word.ActiveDocument

$anchorText = "LOM3081: Introdução à Mecânica dos Sólidos (Requisito fraco)"

$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq $anchorText) {
        $anchorPara = $p
        break
    }
}

$first = $anchorPara.Next()   # empty paragraph right after the anchor
$second = $first.Next()       # "Ver no Jupiter Salvar em pdf Salvar em docx"
$third = $second.Next()       # empty paragraph
$fourth = $third.Next()       # empty paragraph carrying the page break

$start = $first.Range.Start
$end = $fourth.Range.End

$d.Range($start, $end).Delete()

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
